$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(44529, 6213, 11862),
    @(44530, 7432, 14223),
    @(44531, 7856, 14254),
    @(44532, 8600, 16641),
    @(44533, 7352, 17389),
    @(44534, 5636, 13998),
    @(44535, 2711, 6757)
)

$startRow = 286
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r - 1, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

$ws.Range("C285").Select()
